# Add a new "administration_tutelle" column before the existing
# "gestionnaires_additionnels" column (K), pushing that column to L.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at K; everything from K onward (just
# "gestionnaires_additionnels" in K1) shifts one column right to L.
$ws.Columns.Item(11).Insert()

# New header text for the inserted column.
$ws.Range("K1").Value = "administration_tutelle"

# Touching (no-op) a formatting property normalizes the new cell's style
# back to the default/unstyled format without materializing other cells
# in the column.
$ws.Range("K1").WrapText = $false

# Match the new column's width.
$ws.Columns.Item(11).ColumnWidth = 26.83

# Update selection to the new active cell.
$ws.Range("K2").Select()
